# Applies the "same vaccination database" update to Table4:
#  - Swap the order of the "Diabetes" and "Asthma" comorbidity labels
#    (Asthma now appears in row 10, Diabetes in row 11)
#  - Update a handful of Count values to match the refreshed database

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the comorbidity labels in A10/A11
$ws.Range("A10").Value = "Asthma"
$ws.Range("A11").Value = "Diabetes"

# Updated counts
$ws.Range("B5").Value = 34
$ws.Range("B6").Value = 24
$ws.Range("B8").Value = 13
$ws.Range("B10").Value = 5
